$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix column header names (typo/naming correction)
$ws.Range("G1").Value = "Body.damage"
$ws.Range("H1").Value = "Forewing.dorsal.damage"
$ws.Range("I1").Value = "Forewing.ventral.damage"
$ws.Range("J1").Value = "Hindwing.dorsal.damage"
$ws.Range("K1").Value = "Hindwing.ventral.damage"
